# Realestate Update resale numbers 2024-01-05 10:15
# Appends a new data row (row 20) to the CityResaleNum sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Date / Week columns look numeric-ish to Excel's auto-detection (date,
# leading-zero number) so force them to Text first, then restore the
# surrounding cell's (default) style so no stray formatting is left behind.
$ws.Range("A20").NumberFormat = "@"
$ws.Range("A20").Value = "2024-01-05"
$ws.Range("A20").Style = $ws.Range("A19").Style

$ws.Range("B20").Value = "10:15:42"
$ws.Range("C20").Value = "Friday"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "00"
$ws.Range("D20").Style = $ws.Range("D19").Style

$ws.Range("E20").Value = 140500
$ws.Range("F20").Value = 142947
$ws.Range("G20").Value = 171776
$ws.Range("H20").Value = 147028
$ws.Range("I20").Value = -1
$ws.Range("J20").Value = 117749
$ws.Range("K20").Value = 224321
$ws.Range("L20").Value = 248465
$ws.Range("M20").Value = 184725
$ws.Range("N20").Value = 110090
$ws.Range("O20").Value = 40370
$ws.Range("P20").Value = 30792
$ws.Range("Q20").Value = 72382
$ws.Range("R20").Value = -1
$ws.Range("S20").Value = 41208
$ws.Range("T20").Value = -1
